# Delete slide 6 ("Speaking engagement metrics", sldId 322) from the deck.
$p = $ppt.ActivePresentation
$p.Slides.Item(6).Delete()
